$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("C1").Value = "grant_type"

# Update data row 2
$ws.Range("A2").Value = "Prioritise"
$ws.Range("B2").Value = "Password01"
$ws.Range("C2").Value = "password"

# Delete rows 3 through 5 (previously held Johannesburg, Midvaal, Lekwa)
$ws.Range("A3:C5").Delete()
